$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Desmond Bane'
$ws.Range("B2").Value = 'SG,SF'
$ws.Range("C2").Value = 'Memphis Grizzlies'

$ws.Range("A3").Value = 'Scoot Henderson'
$ws.Range("B3").Value = 'PG'
$ws.Range("C3").Value = 'Portland Trail Blazers'

$ws.Range("A4").Value = 'Guerschon Yabusele'
$ws.Range("B4").Value = 'PF,C'
$ws.Range("C4").Value = 'Philadelphia 76ers'

$ws.Range("A5").Value = 'Jamal Murray'
$ws.Range("B5").Value = 'PG,SG'
$ws.Range("C5").Value = 'Denver Nuggets'

$ws.Range("A6").Value = 'Trae Young'
$ws.Range("B6").Value = 'PG'
$ws.Range("C6").Value = 'Atlanta Hawks'

$ws.Range("A7").Value = 'Norman Powell'
$ws.Range("B7").Value = 'SG,SF'
$ws.Range("C7").Value = 'LA Clippers'

$ws.Range("A8").Value = 'Kawhi Leonard'
$ws.Range("B8").Value = 'SG,SF,PF'
$ws.Range("C8").Value = 'LA Clippers'

$ws.Range("A9").Value = 'LeBron James'
$ws.Range("B9").Value = 'SF,PF'
$ws.Range("C9").Value = 'Los Angeles Lakers'

$ws.Range("A10").Value = 'Bam Adebayo'
$ws.Range("B10").Value = 'C'
$ws.Range("C10").Value = 'Miami Heat'

$ws.Range("A11").Value = 'Myles Turner'
$ws.Range("B11").Value = 'C'
$ws.Range("C11").Value = 'Indiana Pacers'

$ws.Range("A12").Value = 'Cason Wallace'
$ws.Range("B12").Value = 'PG,SG'
$ws.Range("C12").Value = 'Oklahoma City Thunder'

$ws.Range("A13").Value = 'Tari Eason'
$ws.Range("B13").Value = 'SF,PF'
$ws.Range("C13").Value = 'Houston Rockets'

$ws.Range("A14").Value = 'Walker Kessler'
$ws.Range("B14").Value = 'C'
$ws.Range("C14").Value = 'Utah Jazz'

$ws.Range("A15").Value = 'Devin Booker'
$ws.Range("B15").Value = 'PG,SG'
$ws.Range("C15").Value = 'Phoenix Suns'

$ws.Range("A16").Value = 'Jalen Brunson'
$ws.Range("B16").Value = 'PG'
$ws.Range("C16").Value = 'New York Knicks'

$ws.Range("A17").Value = 'Immanuel Quickley'
$ws.Range("B17").Value = 'PG,SG'
$ws.Range("C17").Value = 'Toronto Raptors'

$ws.Range("A18").Value = 'Brandon Ingram'
$ws.Range("B18").Value = 'SG,SF,PF'
$ws.Range("C18").Value = 'New Orleans Pelicans'

$ws.Range("A19").Value = 'D''Angelo Russell'
$ws.Range("B19").Value = 'PG'
$ws.Range("C19").Value = 'Brooklyn Nets'
